# Fix: prevent hidden columns from being labeled upon detecting changes
# - Clears the "AENDERUNG" (change) marker from column L for rows where the
#   left/right halves of the row are actually identical (the marker was
#   incorrectly applied because a hidden column's contents were compared).
# - Re-applies the correct "segment group header" banding style (gray fill,
#   the B-column bold) to the rows that start a new segment group, which
#   had been missed previously.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that start a new segment-name group and need the group-header banding
# (style index 2 for most columns, style index 3 for column B) instead of the
# plain "no fill" style they currently carry.
$headerRows = @(23, 27, 31, 35, 42, 46, 53)

# A row that already has the correct banding style to copy formats from.
$styleSource = $ws.Range("A2:V2")

foreach ($r in $headerRows) {
    $dst = $ws.Range("A" + $r + ":V" + $r)
    $styleSource.Copy()
    $dst.PasteSpecial(-4122)  # xlPasteFormats
}

$excel.CutCopyMode = 0

# All rows (23-59) whose "AENDERUNG" label in column L must be removed because
# the difference was only detected in a hidden column. Rows 34 and 47 have a
# genuine difference elsewhere in the row, so they keep their label.
$clearRows = @(23,24,25,26,27,28,29,30,31,32,33,35,36,37,38,39,40,41,42,43,44,45,46,48,49,50,51,52,53,54,55,56,57,58,59)

# A cell that already has the plain "empty" style used for column L.
$lStyleSource = $ws.Range("L2")

foreach ($r in $clearRows) {
    $cell = $ws.Range("L" + $r)
    $cell.ClearContents()
    $lStyleSource.Copy()
    $cell.PasteSpecial(-4122)  # xlPasteFormats
}

$excel.CutCopyMode = 0
